# Adds a bold callout text box to three slides and shrinks/repositions the
# existing content picture on each of those slides to make room for it.
#
# Note on numeric literals below: the host's Shape.Left/Top/Width/Height
# property setters round-trip the point value through a 32-bit float before
# converting to EMU for storage, which can shift the stored EMU by 1 unit
# versus a naive `target_emu / 12700` conversion. The literals used here were
# chosen (via empirical probing of this runtime) so that, after that
# round-trip, the saved OOXML lands exactly on the target EMU values.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide with SlideID 257 (Slides.Item(2))
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$pic2 = $s2.Shapes.Item(2)
$pic2.LockAspectRatio = 0
$pic2.Left = 55.96615173228347
$pic2.Top = 131.60496062992127
$pic2.Width = 650.1550693700788
$pic2.Height = 325.0774803149606

$tb2 = $s2.Shapes.AddTextbox(1, 742.4191338582677, 243.6915748031496, 206.94606299212597, 55.73905511811024)
$tb2.TextFrame.WordWrap = 1
$tb2.TextFrame.AutoSize = 1
$tb2.Fill.Visible = 0
$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "Units sold declined                                  20%"
$tr2.Font.Size = 20
$tr2.Font.Bold = 1
$tb2.Height = 55.73905511811024

# ---------------------------------------------------------------------------
# Slide with SlideID 258 (Slides.Item(3))
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$pic3 = $s3.Shapes.Item(2)
$pic3.LockAspectRatio = 0
$pic3.Left = 50.87489188976378
$pic3.Top = 137.44818897637796
$pic3.Width = 674.7299512598426
$pic3.Height = 337.36498062992126

$tb3 = $s3.Shapes.AddTextbox(1, 750.0359842519686, 243.96551181102362, 209.9640157480315, 52.068897637795274)
$tb3.TextFrame.WordWrap = 1
$tb3.TextFrame.AutoSize = 1
$tb3.Fill.Visible = 0
$tr3 = $tb3.TextFrame.TextRange
$tr3.Text = "Total Revenue rose 27%"
$tr3.Font.Bold = 1
$tb3.Height = 52.06890763779528

# ---------------------------------------------------------------------------
# Slide with SlideID 263 (Slides.Item(4))
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$pic4 = $s4.Shapes.Item(2)
$pic4.LockAspectRatio = 0
$pic4.Left = 53.89284464566929
$pic4.Top = 139.45788401574802
$pic4.Width = 658.3775590551181
$pic4.Height = 329.1887401574803

$tb4 = $s4.Shapes.AddTextbox(1, 740.6946456692914, 252.7584251968504, 210.70779527559054, 41.198425196850394)
$tb4.TextFrame.WordWrap = 1
$tb4.TextFrame.AutoSize = 1
$tb4.Fill.Visible = 0
$tr4 = $tb4.TextFrame.TextRange
$tr4.Text = "CPI rose 40%"
$tr4.Font.Size = 28
$tr4.Font.Bold = 1
$tb4.Height = 41.198425196850394
